$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the name fields
$ws.Range("A2").Value = "Kollapudi"
$ws.Range("B2").Value = "Venu"

# Update the Email cell value and repoint its existing hyperlink
$ws.Range("C2").Value = "venukollapudi@gmail.com"
$ws.Range("C2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:venukollapudi@gmail.com")
$ws.Range("C2").Style = "Hyperlink"

# Update password / confirmPassword and turn them into hyperlinks too
$ws.Range("D2").Value = "Venu@12345"
$ws.Hyperlinks.Add($ws.Range("D2"), "mailto:Venu@12345")
$ws.Range("D2").Style = "Hyperlink"

$ws.Range("E2").Value = "Venu@12345"
$ws.Hyperlinks.Add($ws.Range("E2"), "mailto:Venu@12345")
$ws.Range("E2").Style = "Hyperlink"

# Move the active selection to G2
$ws.Range("G2").Select()
